# Update the "datetimeFigureOut" date placeholder text on the slide master
# and every slide layout from 10/04/2022 to 11/04/2022.

$p = $ppt.ActivePresentation

$oldDate = "10/04/2022"
$newDate = "11/04/2022"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}
